$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Preserve the formatting "donors" we need before we start overwriting the
#    cells that currently carry them. Formats are copied (via Copy +
#    PasteSpecial formats-only) into a far-away scratch column (AA) so the
#    order of the later writes doesn't matter.
#       AA1 <- A4   (bold 18 + border : header style)
#       AA2 <- A5   (plain + border   : body style)
#       AA3 <- D5   (date + border    : date style)
#       AA4 <- D15  (bold 16 + border : big-label style)
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("AA1").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("AA2").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("AA3").PasteSpecial(-4122)

$ws.Range("D15").Copy()
$ws.Range("AA4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 1) Wipe the old layout clean (old dimension was A1:E18) except for the
#    title cell A1, which is untouched. This removes stale column-E data and
#    the old row 15/18 special cells so nothing lingers.
# ---------------------------------------------------------------------------
$ws.Range("A4:E18").Clear()

# ---------------------------------------------------------------------------
# 2) Rebuild the header row (row 4)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 23

$ws.Range("AA1").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("A4").Value = "TAREA"
$ws.Range("B4").Value = "DÍA"
$ws.Range("C4").Value = "DURACION (min)"
$ws.Range("D4").Value = "NOTAS"

# ---------------------------------------------------------------------------
# 3) Body rows 5-24: plain body style everywhere first ...
# ---------------------------------------------------------------------------
$ws.Range("AA2").Copy()
$ws.Range("A5:D24").PasteSpecial(-4122)

# ... then the date style on the whole B5:B14 date column
$ws.Range("AA3").Copy()
$ws.Range("B5:B14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Task / date / duration / notes data, rows 5-14
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Configuración y puesta a punto de todo el entorno de desarrollo (NetBeans, mysql-connector y servidor mysql"
$ws.Range("B5").Value = 40977
$ws.Range("C5").Value = 60

$ws.Range("A6").Value = "Configuración en el entorno de desarrollo del repositorio GIT"
$ws.Range("B6").Value = 40978
$ws.Range("C6").Value = 30

$ws.Range("A7").Value = " Estudio de la construcción de interfaces con el plugin Matisse de NetBeans que trabaja con awt y swing, además de la construcción de una interfaz ejemplo"
$ws.Range("B7").Value = 40978
$ws.Range("C7").Value = 120

$ws.Range("A8").Value = "Desarrollo de un pequeño ejemplo que trabaje con una BD usando el driver JDBC"
$ws.Range("B8").Value = 40981
$ws.Range("C8").Value = 90

$ws.Range("A9").Value = "Comprensión de la arquitectura MVC"
$ws.Range("B9").Value = 40983
$ws.Range("C9").Value = 30

$ws.Range("A10").Value = "Comprensión de la arquitectura diseñada por el equipo de implementación de todo el sistema, especialmente de la diseñada para la gestión de la interfaz y sus eventos"
$ws.Range("B10").Value = 40989
$ws.Range("C10").Value = 120

$ws.Range("A11").Value = "Modelado del panel Datos Voluntario"
$ws.Range("B11").Value = 40990
$ws.Range("C11").Value = 180
$ws.Range("D11").Value = "Tanto tiempo es debido a la primera toma de contacto con el constructor de interfaces, además para adaptarme al estilo definido en los bocetos de interfaz"

$ws.Range("A12").Value = "Modelado del panel Contabilidad"
$ws.Range("B12").Value = 40990
$ws.Range("C12").Value = 60

$ws.Range("A13").Value = "Implementacion de la clase ControladorContabilidad, además de la reestructuración de las clases que enlazarían con la vista de Contabilidad. Se han ordenado los listeners de la vista principal de forma que ahora se sabe a qué panel pertenece cada listener"
$ws.Range("B13").Value = 40990
$ws.Range("C13").Value = 90

$ws.Range("A14").Value = "Implementacion de la clase ControladorDatosVoluntario"
$ws.Range("B14").Value = 40990
$ws.Range("C14").Value = 45

# ---------------------------------------------------------------------------
# 5) Totals block: row 25 (label + sum), row 29 (hours), row 30 (footer link)
# ---------------------------------------------------------------------------
$ws.Range("AA4").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = 20
$ws.Range("B25").Value = "TIEMPO TOTAL DE TRABAJO"

$ws.Range("AA2").Copy()
$ws.Range("C25:D25").PasteSpecial(-4122)
$ws.Range("C25").Formula = "=SUM(C5:C21)"

$excel.CutCopyMode = 0

$ws.Range("B29").Value = "Horas de trabajo"
$ws.Range("C29").Formula = "=C25/60"

$ws.Range("A30").Value = "El tutorial para la construccion de interfaces se encuentra en: http://netbeans.org/kb/docs/java/quickstart-gui.html#getting_started"

# ---------------------------------------------------------------------------
# 6) Clean up the scratch donor column
# ---------------------------------------------------------------------------
$ws.Range("AA1:AA4").Clear()

# ---------------------------------------------------------------------------
# 7) Column widths (stored xlsx width = ColumnWidth + 5/6 on this engine)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 138.498697916667
$ws.Columns.Item(2).ColumnWidth = 33.998697916667
$ws.Columns.Item(4).ColumnWidth = 130.166666666667
$ws.Columns.Item(6).ColumnWidth = 127.166666666667

# ---------------------------------------------------------------------------
# 8) View: selection + scroll
# ---------------------------------------------------------------------------
$ws.Range("B15").Select()

Write-Output "edit complete"
